$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "EstudioParcial" (time blocked off to study for upcoming midterms) is a brand
# new entry: new shared string + new purple fill (RGB 112,48,160 = FF7030A0).
$purple = 10498160   # RGB(112, 48, 160) packed as BGR for OLE_COLOR
$green  = 5287936    # same green already used by the "TrabajoExtra" cells (G8/G9/A32)
$xlLineStyleNone = -4142

function Set-EstudioParcial($addr) {
    $rng = $ws.Range($addr)
    # Clear the border first, then the fill, then the text -- this ordering
    # lets the engine collapse straight onto a single clean cell style
    # instead of leaving a transient one behind.
    $rng.Borders.LineStyle = $xlLineStyleNone
    $rng.Interior.Color = $purple
    $rng.Value = "EstudioParcial"
}

function Set-TrabajoExtra($addr) {
    $rng = $ws.Range($addr)
    $rng.Borders.LineStyle = $xlLineStyleNone
    $rng.Interior.Color = $green
    $rng.Value = "TrabajoExtra"
}

# Cells that used to be "Descanso" and become "EstudioParcial"
Set-EstudioParcial "G4"
Set-EstudioParcial "G5"
Set-EstudioParcial "B13"
Set-EstudioParcial "B14"

# New cell added next to the legend's existing "TrabajoExtra" entry (A32)
Set-EstudioParcial "B32"

# A couple of cells that used to be "Descanso" become "TrabajoExtra" instead
Set-TrabajoExtra "B17"
Set-TrabajoExtra "B18"

# Move the active selection like the author left it
$ws.Range("H14").Select()
